{"js": "// Replace the 25 \"A\u00d7B=C\" multiplication-fact answers in the body with\n// their updated values (see commit \"Update master to output generated\n// at 9a8706d\"). Each original string is unique in the document, so a\n// case-sensitive whole-text search safely targets exactly one run.\nconst replacements = [\n  [\"402\u00d74=1608\", \"126\u00d75=630\"],\n  [\"347\u00d72=694\", \"776\u00d72=1552\"],\n  [\"800\u00d76=4800\", \"527\u00d73=1581\"],\n  [\"295\u00d74=1180\", \"953\u00d79=8577\"],\n  [\"661\u00d74=2644\", \"528\u00d73=1584\"],\n  [\"308\u00d72=616\", \"748\u00d73=2244\"],\n  [\"921\u00d75=4605\", \"756\u00d77=5292\"],\n  [\"923\u00d79=8307\", \"610\u00d77=4270\"],\n  [\"543\u00d75=2715\", \"532\u00d74=2128\"],\n  [\"680\u00d75=3400\", \"499\u00d78=3992\"],\n  [\"914\u00d77=6398\", \"303\u00d75=1515\"],\n  [\"242\u00d78=1936\", \"289\u00d76=1734\"],\n  [\"512\u00d72=1024\", \"370\u00d74=1480\"],\n  [\"902\u00d72=1804\", \"525\u00d72=1050\"],\n  [\"618\u00d74=2472\", \"833\u00d78=6664\"],\n  [\"355\u00d72=710\", \"272\u00d74=1088\"],\n  [\"787\u00d74=3148\", \"406\u00d72=812\"],\n  [\"438\u00d79=3942\", \"659\u00d73=1977\"],\n  [\"291\u00d73=873\", \"507\u00d79=4563\"],\n  [\"624\u00d76=3744\", \"298\u00d75=1490\"],\n  [\"714\u00d75=3570\", \"323\u00d77=2261\"],\n  [\"163\u00d78=1304\", \"624\u00d78=4992\"],\n  [\"613\u00d73=1839\", \"198\u00d76=1188\"],\n  [\"692\u00d72=1384\", \"655\u00d77=4585\"],\n  [\"445\u00d78=3560\", \"518\u00d72=1036\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"A\u00d7B=C\" multiplication-fact answers in the document with\n# their updated values (see commit \"Update master to output generated at\n# 9a8706d\"). Each original string is unique in the document, so a plain\n# Find/Replace (wdReplaceAll, but only ever one hit) safely retargets\n# exactly the intended cell.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"402\u00d74=1608\", \"126\u00d75=630\"),\n    @(\"347\u00d72=694\", \"776\u00d72=1552\"),\n    @(\"800\u00d76=4800\", \"527\u00d73=1581\"),\n    @(\"295\u00d74=1180\", \"953\u00d79=8577\"),\n    @(\"661\u00d74=2644\", \"528\u00d73=1584\"),\n    @(\"308\u00d72=616\", \"748\u00d73=2244\"),\n    @(\"921\u00d75=4605\", \"756\u00d77=5292\"),\n    @(\"923\u00d79=8307\", \"610\u00d77=4270\"),\n    @(\"543\u00d75=2715\", \"532\u00d74=2128\"),\n    @(\"680\u00d75=3400\", \"499\u00d78=3992\"),\n    @(\"914\u00d77=6398\", \"303\u00d75=1515\"),\n    @(\"242\u00d78=1936\", \"289\u00d76=1734\"),\n    @(\"512\u00d72=1024\", \"370\u00d74=1480\"),\n    @(\"902\u00d72=1804\", \"525\u00d72=1050\"),\n    @(\"618\u00d74=2472\", \"833\u00d78=6664\"),\n    @(\"355\u00d72=710\", \"272\u00d74=1088\"),\n    @(\"787\u00d74=3148\", \"406\u00d72=812\"),\n    @(\"438\u00d79=3942\", \"659\u00d73=1977\"),\n    @(\"291\u00d73=873\", \"507\u00d79=4563\"),\n    @(\"624\u00d76=3744\", \"298\u00d75=1490\"),\n    @(\"714\u00d75=3570\", \"323\u00d77=2261\"),\n    @(\"163\u00d78=1304\", \"624\u00d78=4992\"),\n    @(\"613\u00d73=1839\", \"198\u00d76=1188\"),\n    @(\"692\u00d72=1384\", \"655\u00d77=4585\"),\n    @(\"445\u00d78=3560\", \"518\u00d72=1036\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2) | Out-Null\n}\n"}
